$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42: extend the "Survey Dao Tests" note, add hours, grow the row ---
$ws.Range("D42").Value = "Indie Project: Brought Survey Dao Tests to 5/5 passing; Noted programming decisions that are needed for searching surveys.  Revised tables (and other files as needed) to structure the roles table as neede for authentication.  Created data for the application database.`nWeek 7: Started videos, following along in project."
$ws.Range("B42").Value = 6.5
$ws.Rows(42).RowHeight = 60

# --- Row 44: replace the old timestamp note with the shorter "Thurs" note ---
$ws.Range("D44").Value = "Thurs -  a little more time than listed."

# --- Move the active selection to D48 (single cell) ---
$ws.Range("D48").Select() | Out-Null
